$d = $word.ActiveDocument

# --- Hunk 1 & 2: merge the split "region" runs into single runs ---
$r1 = $d.Content
$r1.Find.ClearFormatting()
$r1.Find.Execute(", WHO South-East Asia Region ", $false, $false, $false, $false, $false, $true, 1, $false, ", WHO South-East Asia Region ", 2) | Out-Null

$r2 = $d.Content
$r2.Find.ClearFormatting()
$r2.Find.Execute(", and WHO Western Pacific Region ", $false, $false, $false, $false, $false, $true, 1, $false, ", and WHO Western Pacific Region ", 2) | Out-Null

# --- Hunk 3: append page break + "References" heading + Mendeley bibliography SDT ---
$refsXml = @'
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:br w:type="page"/>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>References</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
</w:p>
<w:sdt>
  <w:sdtPr>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:tag w:val="MENDELEY_BIBLIOGRAPHY"/>
    <w:id w:val="-2012589588"/>
    <w:placeholder>
      <w:docPart w:val="DefaultPlaceholder_-1854013440"/>
    </w:placeholder>
  </w:sdtPr>
  <w:sdtContent>
    <w:p>
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:ind w:hanging="640"/>
        <w:divId w:val="688064602"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:kern w:val="0"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>[1]</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:tab/>
        <w:t xml:space="preserve">&#8216;WHO EMRO | Outbreaks | Epidemic and pandemic </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>diseases&#8217;</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>. Accessed: May 24, 2024. [Online]. Available: https://www.emro.who.int/pandemic-epidemic-diseases/outbreaks/index.html</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:ind w:hanging="640"/>
        <w:divId w:val="72897699"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>[2]</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:tab/>
        <w:t>&#8216;Dengue- Global situation&#8217;. Accessed: May 25, 2024. [Online]. Available: https://www.who.int/emergencies/disease-outbreak-news/item/2023-DON498</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:ind w:hanging="640"/>
        <w:divId w:val="2061322595"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>[3]</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:tab/>
        <w:t>&#8216;Weekly bulletins on outbreaks and other emergencies | WHO | Regional Office for Africa&#8217;. Accessed: May 24, 2024. [Online]. Available: https://www.afro.who.int/health-topics/disease-outbreaks/outbreaks-and-other-emergencies-updates?page=0</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:ind w:hanging="640"/>
        <w:divId w:val="104035936"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>[4]</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:tab/>
        <w:t xml:space="preserve">&#8216;PAHO/WHO Data - National Dengue fever </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>cases&#8217;</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>. Accessed: May 24, 2024. [Online]. Available: https://www3.paho.org/data/index.php/en/mnu-topics/indicadores-dengue-en/dengue-nacional-en/252-dengue-pais-ano-en.html</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:ind w:hanging="640"/>
        <w:divId w:val="1303468031"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>[5]</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:tab/>
        <w:t>&#8216;SEAR Epidemiological Bulletins&#8217;. Accessed: May 24, 2024. [Online]. Available: https://www.who.int/southeastasia/outbreaks-and-emergencies/health-emergency-information-risk-assessment/sear-epi-bulletins</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:ind w:hanging="640"/>
        <w:divId w:val="1001196179"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>[6]</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:tab/>
        <w:t>&#8216;Dengue Situation Updates 2023&#8217;. Accessed: May 24, 2024. [Online]. Available: https://iris.who.int/handle/10665/365676</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:ind w:hanging="640"/>
        <w:divId w:val="1836527191"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>[7]</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:tab/>
        <w:t>&#8216;List of WHO regions - Wikipedia&#8217;. Accessed: May 26, 2024. [Online]. Available: https://en.wikipedia.org/wiki/List_of_WHO_regions</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:ind w:hanging="640"/>
        <w:divId w:val="821506442"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>[8]</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:tab/>
        <w:t>&#8216;List of Countries by Continent 2024&#8217;. Accessed: May 26, 2024. [Online]. Available: https://worldpopulationreview.com/country-rankings/list-of-countries-by-continent</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t> </w:t>
      </w:r>
    </w:p>
  </w:sdtContent>
</w:sdt>

'@

$end = $d.Content
$end.Collapse(0)
$end.InsertXML($refsXml)

Write-Output "done"
